$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: response TRUE -> FALSE, confidence 0.97 -> 0.9
$ws.Range("B2").Value = "'FALSE"
$ws.Range("C2").Value = 0.9

# Row 3: response FALSE -> TRUE, confidence 0.7 -> 0.8
$ws.Range("B3").Value = "'TRUE"
$ws.Range("C3").Value = 0.8

# Row 4: confidence 0.85 -> 0.9
$ws.Range("C4").Value = 0.9

# Row 5: confidence 0.95 -> 1
$ws.Range("C5").Value = 1

# Row 6: confidence 0.99 -> 1
$ws.Range("C6").Value = 1

# Row 9: response TRUE -> FALSE, confidence 0.8 -> 1
$ws.Range("B9").Value = "'FALSE"
$ws.Range("C9").Value = 1

# Row 10: confidence 0.9 -> 1
$ws.Range("C10").Value = 1
